$d = $word.ActiveDocument

# The version/date line on the cover page reads:
#   "Version 11.07.05, 2015-12-09"
# and must become:
#   "Version 11.08.01, 2016-02-15"
# Anchor the replacement on the part of the string that actually changes
# ("7.05, 2015-12-09" -> "8.01, 2016-02-15") so the untouched leading runs
# ("Version ", "1", "1", ".", "0") are left completely alone.
$d.Content.Find.Execute("7.05, 2015-12-09", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8.01, 2016-02-15", 2)
